$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated odds values per the diff (cell -> new value)

# Row 2
$ws.Range("Q2").Value = 2.15
$ws.Range("R2").Value = 1.62

# Row 3
$ws.Range("Q3").Value = 2.25
$ws.Range("R3").Value = 1.57

# Row 5
$ws.Range("M5").Value = 1.11
$ws.Range("O5").Value = 1.5
$ws.Range("R5").Value = 1.44
$ws.Range("V5").Value = 1.67

# Row 7
$ws.Range("G7").Value = 2.6
$ws.Range("I7").Value = 2.6
$ws.Range("J7").Value = 3.25
$ws.Range("K7").Value = 2.2
$ws.Range("L7").Value = 3.25
$ws.Range("M7").Value = 1.05
$ws.Range("O7").Value = 1.29
$ws.Range("Q7").Value = 1.95
$ws.Range("R7").Value = 1.9
$ws.Range("U7").Value = 1.75
$ws.Range("W7").Value = 9
$ws.Range("X7").Value = 13
$ws.Range("Z7").Value = 26
$ws.Range("AK7").Value = 26
$ws.Range("AN7").Value = 4.75
$ws.Range("AQ7").Value = 51

# Row 8
$ws.Range("L8").Value = 2.37

# Row 9
$ws.Range("Q9").Value = 1.73
$ws.Range("R9").Value = 2.08

# Row 11
$ws.Range("G11").Value = 1.62

# Row 12
$ws.Range("G12").Value = 2.1

# Row 13
$ws.Range("I13").Value = 1.57

# Row 14
$ws.Range("G14").Value = 4.2
$ws.Range("H14").Value = 4.5
$ws.Range("I14").Value = 1.65
$ws.Range("J14").Value = 4.33
$ws.Range("L14").Value = 2.1
$ws.Range("N14").Value = 26
$ws.Range("Z14").Value = 51
$ws.Range("AA14").Value = 29
$ws.Range("AK14").Value = 15
$ws.Range("AN14").Value = 7
$ws.Range("AO14").Value = 21
$ws.Range("AZ14").Value = 21

# Row 15
$ws.Range("M15").Value = 1.03
$ws.Range("O15").Value = 1.2

# Row 16
$ws.Range("I16").Value = 5.5
$ws.Range("K16").Value = 2.38
$ws.Range("M16").Value = 1.04
$ws.Range("P16").Value = 4
$ws.Range("S16").Value = 1.33
$ws.Range("T16").Value = 3.25
$ws.Range("U16").Value = 1.83
$ws.Range("V16").Value = 1.83
$ws.Range("AA16").Value = 12
$ws.Range("AD16").Value = 8.5
$ws.Range("AT16").Value = 3.25
$ws.Range("AX16").Value = 29

# Row 17
$ws.Range("M17").Value = 1.02
$ws.Range("N17").Value = 21
$ws.Range("O17").Value = 1.11

# Row 18
$ws.Range("G18").Value = 2.1
$ws.Range("I18").Value = 3.3
$ws.Range("J18").Value = 2.75
$ws.Range("M18").Value = 1.05
$ws.Range("O18").Value = 1.29
$ws.Range("S18").Value = 1.4
$ws.Range("T18").Value = 2.75
$ws.Range("W18").Value = 8
$ws.Range("X18").Value = 10
$ws.Range("Z18").Value = 19
$ws.Range("AT18").Value = 2.75
$ws.Range("AU18").Value = 8
$ws.Range("AW18").Value = 5.5
$ws.Range("AX18").Value = 19
$ws.Range("BA18").Value = 81

# Row 19
$ws.Range("M19").Value = 1.03
$ws.Range("O19").Value = 1.2
$ws.Range("Q19").Value = 1.7
$ws.Range("R19").Value = 2.1

# Row 21
$ws.Range("Q21").Value = 1.33
$ws.Range("R21").Value = 3.4

# Row 22
$ws.Range("I22").Value = 2.62

# Row 25
$ws.Range("Q25").Value = 1.93
$ws.Range("R25").Value = 1.93

# Row 26
$ws.Range("G26").Value = 2.38
$ws.Range("I26").Value = 2.8
$ws.Range("J26").Value = 3
$ws.Range("L26").Value = 3.25
$ws.Range("S26").Value = 1.33
$ws.Range("T26").Value = 3.25
$ws.Range("W26").Value = 10
$ws.Range("X26").Value = 13
$ws.Range("Z26").Value = 23
$ws.Range("AA26").Value = 19
$ws.Range("AC26").Value = 13
$ws.Range("AL26").Value = 21
$ws.Range("AM26").Value = 26
$ws.Range("AO26").Value = 13
$ws.Range("AT26").Value = 3.25
$ws.Range("AV26").Value = 41
$ws.Range("AY26").Value = 21
$ws.Range("BA26").Value = 51
$ws.Range("BB26").Value = 126

# Row 28
$ws.Range("G28").Value = 2
$ws.Range("I28").Value = 3.8
$ws.Range("J28").Value = 2.62
$ws.Range("L28").Value = 4
$ws.Range("M28").Value = 1.04
$ws.Range("N28").Value = 10
$ws.Range("O28").Value = 1.25
$ws.Range("U28").Value = 1.73
$ws.Range("V28").Value = 2
$ws.Range("X28").Value = 9.5
$ws.Range("Y28").Value = 9
$ws.Range("AA28").Value = 17
$ws.Range("AE28").Value = 13
$ws.Range("AF28").Value = 41
$ws.Range("AH28").Value = 11
$ws.Range("AI28").Value = 19
$ws.Range("AJ28").Value = 13
$ws.Range("AL28").Value = 29
$ws.Range("AM28").Value = 34
$ws.Range("AO28").Value = 11
$ws.Range("AQ28").Value = 41
$ws.Range("AW28").Value = 5.5

# Row 30
$ws.Range("M30").Value = 1.05
$ws.Range("O30").Value = 1.27

# Row 32
$ws.Range("G32").Value = 2.15
$ws.Range("J32").Value = 2.75
$ws.Range("Q32").Value = 1.8
$ws.Range("R32").Value = 2
$ws.Range("W32").Value = 9
$ws.Range("AP32").Value = 19

# Row 33
$ws.Range("Q33").Value = 2.3
$ws.Range("R33").Value = 1.6

# Row 34
$ws.Range("Q34").Value = 1.85
$ws.Range("R34").Value = 2

# Row 35
$ws.Range("O35").Value = 1.4
$ws.Range("P35").Value = 2.75

# Row 38
$ws.Range("R38").Value = 1.53

# Row 39
$ws.Range("Q39").Value = 1.95
$ws.Range("R39").Value = 1.9
